$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook tab / date stamp was rolled forward from 02-12-2025 to 05-12-2025.
# Renaming the sheet also keeps the "Fodterapisystemer" defined name (which points
# at 'Opdateret d. 02-12-2025'!$A$1:$I$17) in sync automatically.
$ws.Name = "Opdateret d. 05-12-2025"
